$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Row, NewB, NewC, NewD, NewE ($null = leave unchanged)
$updates = @(
    @{Row=2; B=$null; C=$null; D='96.590.35'; E='  -1.68%  '},
    @{Row=3; B=$null; C=$null; D='3.665.17'; E='  +0.96%  '},
    @{Row=4; B=$null; C=$null; D=$null; E='  -0.27%  '},
    @{Row=5; B=$null; C=$null; D='239.82'; E='  -2.38%  '},
    @{Row=6; B=$null; C=$null; D='1.85'; E='  +6.20%  '},
    @{Row=7; B=$null; C=$null; D='656.06'; E='  -0.62%  '},
    @{Row=8; B=$null; C=$null; D='0.420'; E='  -0.08%  '},
    @{Row=9; B=$null; C=$null; D='1.08'; E='  -0.20%  '},
    @{Row=10; B=$null; C=$null; D=$null; E='  +0.09%  '},
    @{Row=11; B=$null; C=$null; D='3.661.29'; E='  +0.93%  '},
    @{Row=12; B=$null; C=$null; D='45.83'; E='  +3.49%  '},
    @{Row=13; B=$null; C=$null; D='0.205'; E='  -0.80%  '},
    @{Row=14; B=$null; C=$null; D='6.77'; E='  +4.17%  '},
    @{Row=15; B=$null; C=$null; D='4.345.09'; E='  +0.92%  '},
    @{Row=16; B=$null; C=$null; D='0.0000266'; E='  +1.57%  '},
    @{Row=17; B=$null; C=$null; D='96.272.96'; E='  -1.86%  '},
    @{Row=18; B=$null; C=$null; D='8.07'; E='  -2.23%  '},
    @{Row=19; B=$null; C=$null; D='3.645.28'; E='  +0.62%  '},
    @{Row=20; B=$null; C=$null; D='12.77'; E='  -1.20%  '},
    @{Row=21; B=$null; C=$null; D='18.62'; E='  +1.67%  '},
    @{Row=22; B=$null; C=$null; D='0.524'; E='  -3.68%  '},
    @{Row=23; B=$null; C=$null; D='522.83'; E='  +1.11%  '},
    @{Row=24; B=$null; C=$null; D='3.43'; E='  -0.93%  '},
    @{Row=25; B=$null; C=$null; D='7.06'; E='  +1.75%  '},
    @{Row=26; B=$null; C=$null; D=$null; E='  -2.71%  '},
    @{Row=27; B=$null; C=$null; D='101.49'; E='  +2.05%  '},
    @{Row=28; B=$null; C=$null; D='13.15'; E='  -0.12%  '},
    @{Row=29; B=$null; C=$null; D='3.861.27'; E='  +0.98%  '},
    @{Row=30; B=$null; C=$null; D=$null; E='  +8.71%  '},
    @{Row=31; B=$null; C=$null; D='12.41'; E='  +4.16%  '},
    @{Row=32; B=$null; C=$null; D=$null; E='  -1.38%  '},
    @{Row=33; B=$null; C=$null; D='0.999'; E='  -0.29%  '},
    @{Row=34; B=$null; C=$null; D='1.89'; E='  +14.30%  '},
    @{Row=35; B=$null; C=$null; D='0.185'; E='  -0.93%  '},
    @{Row=36; B=$null; C=$null; D=$null; E='  -0.04%  '},
    @{Row=37; B='Bittensor'; C='https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'; D='663.62'; E='  +7.61%  '},
    @{Row=38; B='EthereumClassic'; C='https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; D='32.38'; E='  +1.02%  '},
    @{Row=39; B=$null; C=$null; D=$null; E='  +2.35%  '},
    @{Row=40; B=$null; C=$null; D='8.79'; E='  -0.66%  '},
    @{Row=41; B=$null; C=$null; D='39.94'; E='  +20.27%  '},
    @{Row=42; B=$null; C=$null; D=$null; E='  +3.55%  '},
    @{Row=43; B=$null; C=$null; D='1.97'; E='  -0.84%  '},
    @{Row=44; B=$null; C=$null; D='0.957'; E='  +2.24%  '},
    @{Row=45; B=$null; C=$null; D='6.49'; E='  +6.50%  '},
    @{Row=46; B=$null; C=$null; D=$null; E='  +0.01%  '},
    @{Row=47; B=$null; C=$null; D='0.0466'; E='  +4.21%  '},
    @{Row=48; B=$null; C=$null; D='0.444'; E='  +11.60%  '},
    @{Row=49; B=$null; C=$null; D=$null; E='  -0.87%  '},
    @{Row=50; B=$null; C=$null; D='23.65'; E='  -0.17%  '},
    @{Row=51; B='MantraDAO'; C='https://coinranking.com/coin/cTdD8lD-6+mantradao-om'; D='3.63'; E='  +3.08%  '}
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($null -ne $u.B) { $ws.Range("B$r").Value = $u.B }
    if ($null -ne $u.C) { $ws.Range("C$r").Value = $u.C }
    if ($null -ne $u.D) {
        $cell = $ws.Range("D$r")
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.Style = "Normal"
    }
    if ($null -ne $u.E) { $ws.Range("E$r").Value = $u.E }
}
